# Auto-generated COM-interop script applying the CasosColombia.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric value corrections on existing rows ---
$numericUpdates = @{
    "BT9" = 1
    "I15" = 2
    "L18" = 1
    "L19" = 1
    "BU26" = 20
    "CW27" = 1
    "CM30" = 1
    "CM31" = 1
    "BU37" = 33
    "BU45" = 43
    "I60" = 79
    "I61" = 79
    "AY61" = 228
    "AY62" = 229
    "AY63" = 229
    "AY64" = 417
    "AY65" = 429
    "AY66" = 432
    "AY67" = 525
    "AY68" = 716
    "AY69" = 741
    "AY70" = 869
    "AY71" = 922
    "AY72" = 1001
    "AY73" = 1002
    "AY74" = 1005
    "AY75" = 1142
    "AY76" = 1180
    "AY77" = 1271
    "I78" = 126
    "AY78" = 1344
    "I79" = 126
    "AY79" = 1345
    "I80" = 127
    "AY80" = 1347
    "I81" = 128
    "AY81" = 1409
    "AY82" = 1464
    "AY83" = 1493
    "AY84" = 1682
    "AY85" = 1721
    "AY86" = 1757
    "AY87" = 1791
    "AY88" = 1807
    "AY89" = 1807
    "AY90" = 1810
    "BT90" = 55
    "I91" = 157
    "AY91" = 1827
    "BT91" = 62
    "AY92" = 1827
    "AY93" = 1881
    "AY94" = 1908
    "AY95" = 1920
    "AY96" = 1957
    "AY97" = 1961
    "H98" = 20
    "AY98" = 1985
    "AY99" = 1994
    "AY100" = 2010
    "AY101" = 2027
    "AY102" = 2041
    "AY103" = 2046
    "H104" = 44
    "AY104" = 2066
    "AY105" = 2074
    "BU105" = 105
    "AY106" = 2079
    "BU106" = 107
    "AY107" = 2087
    "BU107" = 108
    "AY108" = 2087
    "BQ108" = 467
    "AY109" = 2099
    "I110" = 231
    "AY110" = 2099
    "I111" = 237
    "AY111" = 2106
    "AY112" = 2107
    "AY113" = 2115
    "AY114" = 2126
    "AY115" = 2135
    "BU115" = 153
    "AY116" = 2139
    "BT116" = 122
    "AY117" = 2156
    "BT117" = 126
    "AY118" = 2178
    "BQ118" = 642
    "BT118" = 126
    "AY119" = 2197
    "BQ119" = 666
    "BT119" = 130
    "AY120" = 2203
    "BQ120" = 669
    "AY121" = 2203
    "BQ121" = 703
    "AY122" = 2210
    "BQ122" = 769
    "I123" = 271
    "AY123" = 2239
    "BQ123" = 792
    "AY124" = 2258
    "BQ124" = 810
    "AY125" = 2264
    "BQ125" = 824
    "AY126" = 2267
    "BQ126" = 839
    "AY127" = 2272
    "BQ127" = 852
    "CL127" = 187
    "AY128" = 2278
    "BQ128" = 876
    "BU128" = 253
    "AY129" = 2281
    "BQ129" = 884
    "BU129" = 260
    "AY130" = 2288
    "BQ130" = 909
    "BU130" = 262
    "AY131" = 2297
    "BQ131" = 924
    "BU131" = 265
    "AY132" = 2303
    "BQ132" = 948
    "AY133" = 2303
    "BQ133" = 987
    "AY134" = 2312
    "BQ134" = 1053
    "BU134" = 332
    "AY135" = 2314
    "BQ135" = 1092
    "BU135" = 338
    "AY136" = 2321
    "BQ136" = 1141
    "BU136" = 351
    "AY137" = 2321
    "BQ137" = 1234
    "BU137" = 355
    "AY138" = 2323
    "BQ138" = 1338
    "BU138" = 377
    "I139" = 416
    "AY139" = 2323
    "BQ139" = 1348
    "I140" = 508
    "AY140" = 2324
    "BQ140" = 1401
    "I141" = 549
    "AY141" = 2328
    "BQ141" = 1470
    "I142" = 581
    "AY142" = 2336
    "BQ142" = 1522
    "AY143" = 2354
    "BQ143" = 1576
    "BT143" = 314
    "BU143" = 497
    "AY144" = 2360
    "BQ144" = 1634
    "BU144" = 516
    "I145" = 722
    "AY145" = 2363
    "BQ145" = 1719
    "BU145" = 535
    "I146" = 733
    "AY146" = 2363
    "BQ146" = 1862
    "BU146" = 560
    "I147" = 792
    "AY147" = 2364
    "BQ147" = 1936
    "BU147" = 561
    "I148" = 849
    "AY148" = 2365
    "BQ148" = 1997
    "BU148" = 591
    "I149" = 956
    "AY149" = 2365
    "BQ149" = 2092
    "I150" = 988
    "AY150" = 2366
    "BQ150" = 2242
    "I151" = 1018
    "AY151" = 2368
    "BQ151" = 2360
    "BT151" = 534
    "CL151" = 576
    "I152" = 1033
    "AY152" = 2368
    "BQ152" = 2392
    "I153" = 1057
    "AY153" = 2368
    "BQ153" = 2533
    "BU153" = 657
    "I154" = 1108
    "AY154" = 2379
    "BQ154" = 2623
    "BU154" = 710
    "I155" = 1157
    "AY155" = 2381
    "BQ155" = 2737
    "I156" = 1191
    "AY156" = 2381
    "BQ156" = 2912
    "I157" = 1216
    "AY157" = 2384
    "BQ157" = 3017
    "I158" = 1313
    "AY158" = 2384
    "BQ158" = 3088
    "BU158" = 828
    "I159" = 1384
    "AY159" = 2392
    "BQ159" = 3124
    "BU159" = 833
    "I160" = 1417
    "AY160" = 2397
    "BQ160" = 3364
    "BU160" = 869
    "I161" = 1456
    "AY161" = 2411
    "BQ161" = 3414
    "BU161" = 921
    "I162" = 1486
    "BQ162" = 3585
    "BU162" = 950
    "AY163" = 2419
    "BQ163" = 3723
    "BU163" = 1006
    "AY164" = 2424
    "BQ164" = 3837
    "BU164" = 1043
    "AY165" = 2431
    "BQ165" = 4063
    "I166" = 1687
    "BQ166" = 4199
    "BT166" = 877
    "I167" = 1736
    "BQ167" = 4301
    "I168" = 1823
    "BQ168" = 4610
}
foreach ($cellRef in $numericUpdates.Keys) {
    $ws.Range($cellRef).Value = $numericUpdates[$cellRef]
}

# --- Cells reverted to "NaN" (missing-data marker) ---
$nanCells = @(
    "BT12"
    "BU12"
    "BT13"
    "BT14"
    "BT15"
    "BT16"
    "BQ18"
    "CM27"
    "CM28"
    "H90"
    "CF92"
    "DN59"
    "L87"
    "CL120"
    "BU132"
    "BU133"
    "H134"
    "CL140"
    "H147"
    "CS152"
)
foreach ($cellRef in $nanCells) {
    $ws.Range($cellRef).Value = "NaN"
}

# --- New row 170 (2020-08-21) appended to the dataset ---
$row170 = @{
    "A170" = 44064
    "B170" = 522138
    "C170" = 2678
    "D170" = 67893
    "E170" = 62464
    "F170" = 181775
    "G170" = 22738
    "H170" = 2513
    "I170" = 2120
    "J170" = 4603
    "K170" = 3837
    "L170" = 7187
    "M170" = 3592
    "N170" = 16727
    "O170" = 18238
    "P170" = 4222
    "Q170" = 3078
    "R170" = 10996
    "S170" = 5489
    "T170" = 12363
    "U170" = 8325
    "V170" = 2328
    "W170" = 803
    "X170" = 4050
    "Y170" = 12343
    "Z170" = 10113
    "AA170" = 5222
    "AB170" = 41887
    "AC170" = 798
    "AD170" = 117
    "AE170" = 169
    "AF170" = 437
    "AG170" = 28
    "AH170" = 15
    "AI170" = 216
    "AJ170" = 1923
    "AK170" = 2158
    "AL170" = 35226
    "AM170" = 5677
    "AN170" = 2357
    "AO170" = 33034
    "AP170" = 794
    "AQ170" = 18708
    "AR170" = 1398
    "AS170" = 5799
    "AT170" = 1371
    "AU170" = 1532
    "AV170" = 2958
    "AW170" = 1390
    "AX170" = 925
    "AY170" = 2446
    "AZ170" = 2559
    "BA170" = 38804
    "BB170" = 10795
    "BC170" = 1735
    "BD170" = 6482
    "BE170" = 2594
    "BF170" = 274
    "BG170" = 1364
    "BH170" = 2503
    "BI170" = 725
    "BJ170" = 1909
    "BK170" = 7434
    "BL170" = 7300
    "BM170" = 6957
    "BN170" = 13475
    "BO170" = 1838
    "BP170" = 739
    "BQ170" = 5199
    "BR170" = 4544
    "BS170" = 5045
    "BT170" = 1114
    "BU170" = 1267
    "BV170" = 1983
    "BW170" = 2301
    "BX170" = 592
    "BY170" = 3694
    "BZ170" = 2110
    "CA170" = 925
    "CB170" = 593
    "CC170" = 1703
    "CD170" = 1714
    "CE170" = 839
    "CF170" = 720
    "CG170" = 3866
    "CH170" = 968
    "CI170" = 1010
    "CJ170" = 1020
    "CK170" = 1331
    "CL170" = 1201
    "CM170" = 1149
    "CN170" = 970
    "CO170" = 932
    "CP170" = 1018
    "CQ170" = 487
    "CR170" = 2798
    "CS170" = 800
    "CT170" = 731
    "CU170" = 639
    "CV170" = 1077
    "CW170" = 965
    "CX170" = 557
    "CY170" = 659
    "CZ170" = 722
    "DA170" = 1019
    "DB170" = 824
    "DC170" = 931
    "DD170" = 743
    "DE170" = 309
    "DF170" = 323
    "DG170" = 622
    "DH170" = 504
    "DI170" = 378
    "DJ170" = 517
    "DK170" = 300
    "DL170" = 513
    "DM170" = 692
    "DN170" = 498
    "DO170" = 469
    "DP170" = 343
    "DQ170" = 508
    "DR170" = 107607
    "DS170" = 218036
    "DT170" = 7246
    "DU170" = 94723
    "DV170" = 61679
    "DW170" = 20717
    "DX170" = 6959
}
foreach ($cellRef in $row170.Keys) {
    $ws.Range($cellRef).Value = $row170[$cellRef]
}
